$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb2"
$ws.Range("C2").Value = "Ephb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 52.61615766666667
$ws.Range("H2").Value = 157.848473
$ws.Range("I2").Value = 0.7671520491359202
$ws.Range("J2").Value = 0.7671520491359202
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.501929333333333
$ws.Range("N2").Value = 4.505788
$ws.Range("O2").Value = 0.7650463650777426
$ws.Range("P2").Value = 0.7650463650777426
$ws.Range("Q2").Value = 79.02575060685822
$ws.Range("R2").Value = 711.231755461724
$ws.Range("S2").Value = 0.5869068866533775
$ws.Range("T2").Value = 0.5869068866533775

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb2"
$ws.Range("C3").Value = "Ephb1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 52.61615766666667
$ws.Range("H3").Value = 157.848473
$ws.Range("I3").Value = 0.7671520491359202
$ws.Range("J3").Value = 0.7671520491359202
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.461258
$ws.Range("N3").Value = 1.383774
$ws.Range("O3").Value = 0.2349536349222574
$ws.Range("P3").Value = 0.2349536349222574
$ws.Range("Q3").Value = 24.26962365301134
$ws.Range("R3").Value = 218.426612877102
$ws.Range("S3").Value = 0.1802451624825427
$ws.Range("T3").Value = 0.1802451624825427

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efnb2"
$ws.Range("C4").Value = "Ephb1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.395935333333332
$ws.Range("H4").Value = 28.187806
$ws.Range("I4").Value = 0.1369942497546098
$ws.Range("J4").Value = 0.1369942497546098
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.501929333333333
$ws.Range("N4").Value = 4.505788
$ws.Range("O4").Value = 0.7650463650777426
$ws.Range("P4").Value = 0.7650463650777426
$ws.Range("Q4").Value = 14.11203089123644
$ws.Range("R4").Value = 127.008278021128
$ws.Range("S4").Value = 0.1048069528113167
$ws.Range("T4").Value = 0.1048069528113167

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb2"
$ws.Range("C5").Value = "Ephb1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.395935333333332
$ws.Range("H5").Value = 28.187806
$ws.Range("I5").Value = 0.1369942497546098
$ws.Range("J5").Value = 0.1369942497546098
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.461258
$ws.Range("N5").Value = 1.383774
$ws.Range("O5").Value = 0.2349536349222574
$ws.Range("P5").Value = 0.2349536349222574
$ws.Range("Q5").Value = 4.333950339982666
$ws.Range("R5").Value = 39.005553059844
$ws.Range("S5").Value = 0.03218729694329314
$ws.Range("T5").Value = 0.03218729694329314

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Efnb2"
$ws.Range("C6").Value = "Ephb1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.08161033333333334
$ws.Range("H6").Value = 0.244831
$ws.Range("I6").Value = 0.001189891797952309
$ws.Range("J6").Value = 0.001189891797952309
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.501929333333333
$ws.Range("N6").Value = 4.505788
$ws.Range("O6").Value = 0.7650463650777426
$ws.Range("P6").Value = 0.7650463650777426
$ws.Range("Q6").Value = 0.1225729535364444
$ws.Range("R6").Value = 1.103156581828
$ws.Range("S6").Value = 0.0009103223948592336
$ws.Range("T6").Value = 0.0009103223948592336

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Efnb2"
$ws.Range("C7").Value = "Ephb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.08161033333333334
$ws.Range("H7").Value = 0.244831
$ws.Range("I7").Value = 0.001189891797952309
$ws.Range("J7").Value = 0.001189891797952309
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.461258
$ws.Range("N7").Value = 1.383774
$ws.Range("O7").Value = 0.2349536349222574
$ws.Range("P7").Value = 0.2349536349222574
$ws.Range("Q7").Value = 0.03764341913266667
$ws.Range("R7").Value = 0.3387907721940001
$ws.Range("S7").Value = 0.0002795694030930752
$ws.Range("T7").Value = 0.0002795694030930752

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efnb2"
$ws.Range("C8").Value = "Ephb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.492645
$ws.Range("H8").Value = 19.477935
$ws.Range("I8").Value = 0.09466380931151776
$ws.Range("J8").Value = 0.09466380931151776
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.501929333333333
$ws.Range("N8").Value = 4.505788
$ws.Range("O8").Value = 0.7650463650777426
$ws.Range("P8").Value = 0.7650463650777426
$ws.Range("Q8").Value = 9.75149397642
$ws.Range("R8").Value = 87.76344578778001
$ws.Range("S8").Value = 0.07242220321818922
$ws.Range("T8").Value = 0.07242220321818922

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efnb2"
$ws.Range("C9").Value = "Ephb1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.492645
$ws.Range("H9").Value = 19.477935
$ws.Range("I9").Value = 0.09466380931151776
$ws.Range("J9").Value = 0.09466380931151776
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.461258
$ws.Range("N9").Value = 1.383774
$ws.Range("O9").Value = 0.2349536349222574
$ws.Range("P9").Value = 0.2349536349222574
$ws.Range("Q9").Value = 2.99478444741
$ws.Range("R9").Value = 26.95306002669
$ws.Range("S9").Value = 0.02224160609332853
$ws.Range("T9").Value = 0.02224160609332853
